# Update the DATA sheet: the 'password' column (F) for the first three
# data rows now stores the base64-encoded credential
# ("admi123" -> "YWRtaTEyMw==") instead of the old plaintext shared
# string, matching the new DecodeUtils.getDecodedString() usage.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Activate()

# Leading apostrophe forces these to stay plain text (preserves the
# existing quote-prefixed cell style instead of Excel reinterpreting it).
$ws.Range("F2").Value = "'YWRtaTEyMw=="
$ws.Range("F3").Value = "'YWRtaTEyMw=="
$ws.Range("F4").Value = "'YWRtaTEyMw=="

# Move/restore the active selection to F2, as in the edited workbook.
$ws.Range("F2").Select()
